$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "B2" = -0.5105544497887244
    "C2" = 0.4677058628456969
    "D2" = -0.4252444871073068
    "E2" = -0.1454514859615633

    "B3" = -0.704984082927168
    "C3" = 1.048539417301345
    "D3" = -0.5115076716610412
    "E3" = 1.128941192698981

    "B4" = -0.3077101733269968
    "C4" = 0.7499702733889624
    "D4" = -0.387488401597478
    "E4" = -1.066945061397698

    "B5" = 0.3194398239056268
    "C5" = -0.5509028080077042
    "D5" = -0.9913314054426073
    "E5" = 0.9453299126560957

    "B6" = -0.1336755871392932
    "C6" = -1.015835298249297
    "D6" = -1.294854530248012
    "E6" = -0.7199519196186355

    "B7" = -1.226069064553532
    "C7" = -0.06946859387738215
    "D7" = -0.2533835077588892
    "E7" = 1.097845061231636
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}
